$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2, 1).Value2 = "Última actualización: 11:11:31"

$ws.Cells.Item(3, 1).Value2 = "Total filas: 184"

$ws.Cells.Item(37, 1).Value2 = "05:42:22"
$ws.Cells.Item(37, 3).Value2 = "17X38_ROMERO"
$ws.Cells.Item(37, 4).Value2 = 114

$ws.Cells.Item(38, 1).Value2 = "06:33:46"
$ws.Cells.Item(38, 3).Value2 = "27_EL RETIRO"
$ws.Cells.Item(38, 4).Value2 = 63

$ws.Cells.Item(50, 1).Value2 = "07:12:53"
$ws.Cells.Item(50, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(50, 4).Value2 = 49

$ws.Cells.Item(51, 1).Value2 = "06:45:50"
$ws.Cells.Item(51, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(51, 4).Value2 = 76

$ws.Cells.Item(84, 1).Value2 = "08:39:08"
$ws.Cells.Item(84, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(84, 4).Value2 = 25

$ws.Cells.Item(85, 1).Value2 = "07:36:59"
$ws.Cells.Item(85, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(85, 4).Value2 = 88

$ws.Cells.Item(107, 1).Value2 = "09:21:49"
$ws.Cells.Item(107, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(107, 4).Value2 = 42

$ws.Cells.Item(108, 1).Value2 = "08:11:27"
$ws.Cells.Item(108, 3).Value2 = "215C_EL PATO"
$ws.Cells.Item(108, 4).Value2 = 112

$ws.Cells.Item(130, 1).Value2 = "08:39:08"
$ws.Cells.Item(130, 3).Value2 = "16_P MOR-SANTA ANA"
$ws.Cells.Item(130, 4).Value2 = 118

$ws.Cells.Item(131, 1).Value2 = "10:36:18"
$ws.Cells.Item(131, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(131, 4).Value2 = 1

$ws.Cells.Item(140, 3).Value2 = "10_OLMOS"

$ws.Cells.Item(141, 3).Value2 = "16_SANTA ANA"

$ws.Cells.Item(149, 1).Value2 = "10:04:17"
$ws.Cells.Item(149, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(149, 4).Value2 = 67

$ws.Cells.Item(150, 1).Value2 = "11:11:31"
$ws.Cells.Item(150, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(150, 4).Value2 = 0

$ws.Cells.Item(151, 1).Value2 = "10:36:18"
$ws.Cells.Item(151, 2).Value2 = "11:11"
$ws.Cells.Item(151, 3).Value2 = "15_ABASTO"
$ws.Cells.Item(151, 4).Value2 = 35

$ws.Cells.Item(152, 2).Value2 = "11:14"
$ws.Cells.Item(152, 3).Value2 = "225_C ROCA-H SUR"
$ws.Cells.Item(152, 4).Value2 = 113

$ws.Cells.Item(153, 1).Value2 = "09:21:49"
$ws.Cells.Item(153, 2).Value2 = "11:20"
$ws.Cells.Item(153, 4).Value2 = 119

$ws.Cells.Item(154, 1).Value2 = "10:04:17"
$ws.Cells.Item(154, 3).Value2 = "215C_EL PATO"
$ws.Cells.Item(154, 4).Value2 = 77

$ws.Cells.Item(155, 1).Value2 = "10:48:14"
$ws.Cells.Item(155, 2).Value2 = "11:21"
$ws.Cells.Item(155, 4).Value2 = 33

$ws.Cells.Item(156, 2).Value2 = "11:22"
$ws.Cells.Item(156, 3).Value2 = "10_OLMOS"
$ws.Cells.Item(156, 4).Value2 = 46

$ws.Cells.Item(157, 2).Value2 = "11:24"
$ws.Cells.Item(157, 3).Value2 = "11_ETCHEVERRY"
$ws.Cells.Item(157, 4).Value2 = 48

$ws.Cells.Item(158, 1).Value2 = "10:36:18"
$ws.Cells.Item(158, 3).Value2 = "16_P MOR-SANTA ANA"
$ws.Cells.Item(158, 4).Value2 = 49

$ws.Cells.Item(159, 2).Value2 = "11:25"
$ws.Cells.Item(159, 3).Value2 = "11_ETCHEVERRY"
$ws.Cells.Item(159, 4).Value2 = 81

$ws.Cells.Item(160, 1).Value2 = "10:04:17"
$ws.Cells.Item(160, 2).Value2 = "11:30"
$ws.Cells.Item(160, 3).Value2 = "15X38_ABASTO"
$ws.Cells.Item(160, 4).Value2 = 86

$ws.Cells.Item(161, 2).Value2 = "11:32"
$ws.Cells.Item(161, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(161, 4).Value2 = 44

$ws.Cells.Item(162, 1).Value2 = "10:48:14"
$ws.Cells.Item(162, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(162, 4).Value2 = 45

$ws.Cells.Item(163, 1).Value2 = "10:36:18"
$ws.Cells.Item(163, 2).Value2 = "11:33"
$ws.Cells.Item(163, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(163, 4).Value2 = 57

$ws.Cells.Item(164, 1).Value2 = "10:04:17"
$ws.Cells.Item(164, 2).Value2 = "11:34"
$ws.Cells.Item(164, 3).Value2 = "10_OLMOS"
$ws.Cells.Item(164, 4).Value2 = 90

$ws.Cells.Item(165, 1).Value2 = "10:36:18"
$ws.Cells.Item(165, 2).Value2 = "11:35"
$ws.Cells.Item(165, 3).Value2 = "16_SANTA ANA"
$ws.Cells.Item(165, 4).Value2 = 59

$ws.Cells.Item(166, 2).Value2 = "11:37"
$ws.Cells.Item(166, 3).Value2 = "16_P MOR-SANTA ANA"
$ws.Cells.Item(166, 4).Value2 = 93

$ws.Cells.Item(167, 1).Value2 = "10:04:17"
$ws.Cells.Item(167, 2).Value2 = "11:40"
$ws.Cells.Item(167, 3).Value2 = "215A_EL PATO"
$ws.Cells.Item(167, 4).Value2 = 96

$ws.Cells.Item(168, 1).Value2 = "10:55:25"
$ws.Cells.Item(168, 2).Value2 = "11:44"
$ws.Cells.Item(168, 4).Value2 = 49

$ws.Cells.Item(169, 1).Value2 = "10:04:17"
$ws.Cells.Item(169, 2).Value2 = "11:45"
$ws.Cells.Item(169, 3).Value2 = "215B_EL PATO"
$ws.Cells.Item(169, 4).Value2 = 101

$ws.Cells.Item(170, 1).Value2 = "10:55:25"
$ws.Cells.Item(170, 2).Value2 = "11:53"
$ws.Cells.Item(170, 3).Value2 = "15_ABASTO"
$ws.Cells.Item(170, 4).Value2 = 58

$ws.Cells.Item(171, 1).Value2 = "10:04:17"
$ws.Cells.Item(171, 2).Value2 = "11:54"
$ws.Cells.Item(171, 3).Value2 = "225_GOMEZ"
$ws.Cells.Item(171, 4).Value2 = 110

$ws.Cells.Item(172, 1).Value2 = "11:11:31"
$ws.Cells.Item(172, 2).Value2 = "12:05"
$ws.Cells.Item(172, 3).Value2 = "17_ROMERO"
$ws.Cells.Item(172, 4).Value2 = 54

$ws.Cells.Item(173, 1).Value2 = "11:11:31"
$ws.Cells.Item(173, 2).Value2 = "12:06"
$ws.Cells.Item(173, 3).Value2 = "23_HERNANDEZ"
$ws.Cells.Item(173, 4).Value2 = 55

$ws.Cells.Item(174, 1).Value2 = "10:48:14"
$ws.Cells.Item(174, 2).Value2 = "12:07"
$ws.Cells.Item(174, 3).Value2 = "14_ABASTO"
$ws.Cells.Item(174, 4).Value2 = 79

$ws.Cells.Item(175, 1).Value2 = "11:11:31"
$ws.Cells.Item(175, 2).Value2 = "12:17"
$ws.Cells.Item(175, 3).Value2 = "15_ABASTO"
$ws.Cells.Item(175, 4).Value2 = 66

$ws.Cells.Item(176, 1).Value2 = "11:11:31"
$ws.Cells.Item(176, 2).Value2 = "12:18"
$ws.Cells.Item(176, 3).Value2 = "10_OLMOS"
$ws.Cells.Item(176, 4).Value2 = 67

$ws.Cells.Item(177, 1).Value2 = "10:36:18"
$ws.Cells.Item(177, 2).Value2 = "12:29"
$ws.Cells.Item(177, 3).Value2 = "215C_EL PATO"
$ws.Cells.Item(177, 4).Value2 = 113

$ws.Cells.Item(178, 1).Value2 = "10:36:18"
$ws.Cells.Item(178, 2).Value2 = "12:30"
$ws.Cells.Item(178, 3).Value2 = "11_ETCHEVERRY"
$ws.Cells.Item(178, 4).Value2 = 114

$ws.Cells.Item(179, 1).Value2 = "10:36:18"
$ws.Cells.Item(179, 2).Value2 = "12:31"
$ws.Cells.Item(179, 3).Value2 = "16_P MOR-SANTA ANA"
$ws.Cells.Item(179, 4).Value2 = 115

$ws.Cells.Item(180, 2).Value2 = "12:31"
$ws.Cells.Item(180, 3).Value2 = "11_ETCHEVERRY"
$ws.Cells.Item(180, 4).Value2 = 103

$ws.Cells.Item(181, 2).Value2 = "12:36"
$ws.Cells.Item(181, 3).Value2 = "27_EL RETIRO"
$ws.Cells.Item(181, 4).Value2 = 101

$ws.Cells.Item(182, 1).Value2 = "10:48:14"
$ws.Cells.Item(182, 2).Value2 = "12:37"
$ws.Cells.Item(182, 3).Value2 = "27_EL RETIRO"
$ws.Cells.Item(182, 4).Value2 = 109
$ws.Cells.Item(182, 5).Value2 = "LP1912"

$ws.Cells.Item(183, 1).Value2 = "10:48:14"
$ws.Cells.Item(183, 2).Value2 = "12:40"
$ws.Cells.Item(183, 3).Value2 = "15X38_ABASTO"
$ws.Cells.Item(183, 4).Value2 = 112
$ws.Cells.Item(183, 5).Value2 = "LP1912"

$ws.Cells.Item(184, 1).Value2 = "10:55:25"
$ws.Cells.Item(184, 2).Value2 = "12:42"
$ws.Cells.Item(184, 3).Value2 = "14_ABASTO"
$ws.Cells.Item(184, 4).Value2 = 107
$ws.Cells.Item(184, 5).Value2 = "LP1912"

$ws.Cells.Item(185, 1).Value2 = "10:55:25"
$ws.Cells.Item(185, 2).Value2 = "12:43"
$ws.Cells.Item(185, 3).Value2 = "15X38_ABASTO"
$ws.Cells.Item(185, 4).Value2 = 108
$ws.Cells.Item(185, 5).Value2 = "LP1912"

$ws.Cells.Item(186, 1).Value2 = "10:48:14"
$ws.Cells.Item(186, 2).Value2 = "12:43"
$ws.Cells.Item(186, 3).Value2 = "14_ABASTO"
$ws.Cells.Item(186, 4).Value2 = 115
$ws.Cells.Item(186, 5).Value2 = "LP1912"

$ws.Cells.Item(187, 1).Value2 = "11:11:31"
$ws.Cells.Item(187, 2).Value2 = "12:54"
$ws.Cells.Item(187, 3).Value2 = "15X38_ABASTO"
$ws.Cells.Item(187, 4).Value2 = 103
$ws.Cells.Item(187, 5).Value2 = "LP1912"

$ws.Cells.Item(188, 1).Value2 = "11:11:31"
$ws.Cells.Item(188, 2).Value2 = "13:01"
$ws.Cells.Item(188, 3).Value2 = "215C_EL PATO"
$ws.Cells.Item(188, 4).Value2 = 110
$ws.Cells.Item(188, 5).Value2 = "LP1912"

$ws.Cells.Item(189, 1).Value2 = "11:11:31"
$ws.Cells.Item(189, 2).Value2 = "13:06"
$ws.Cells.Item(189, 3).Value2 = "14_ABASTO"
$ws.Cells.Item(189, 4).Value2 = 115
$ws.Cells.Item(189, 5).Value2 = "LP1912"


# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2, 1).Value2 = "Última actualización: 11:11:31"

$ws.Cells.Item(3, 1).Value2 = "Total filas: 26"

$ws.Cells.Item(31, 1).Value2 = "11:11:31"
$ws.Cells.Item(31, 2).Value2 = "13:01"
$ws.Cells.Item(31, 3).Value2 = "215C_EL PATO"
$ws.Cells.Item(31, 4).Value2 = 110
$ws.Cells.Item(31, 5).Value2 = "LP1912"


# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2, 1).Value2 = "Última actualización: 11:11:31"

$ws.Cells.Item(3, 1).Value2 = "Total filas: 21"

$ws.Cells.Item(26, 1).Value2 = "11:11:31"
$ws.Cells.Item(26, 2).Value2 = "13:09"
$ws.Cells.Item(26, 3).Value2 = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(26, 4).Value2 = 118
$ws.Cells.Item(26, 5).Value2 = "L6173"

